# constraints-ph1.docx: "add one more constraint; small changes to ER"
#
# Adds a new blank paragraph followed by an "Assumptions:" section (two
# assumption bullets) after the existing "Logical Constraints:" list, and
# appends a final constraint ("Project names are unique") onto the
# paragraph that carries the trailing _GoBack bookmark.

$d = $word.ActiveDocument

# 1) Insert a blank paragraph right after the last existing constraint
#    ("...previously rejected from"), mirroring pressing Enter at the end
#    of that line.
$d.Content.Find.Execute(
    "previously rejected from", $true, $false, $false, $false, $false,
    $true, 1, $false, "previously rejected from^p", 2)

# 2) The trailing paragraph (holding the _GoBack bookmark) is now the
#    last paragraph in the document. Type the new "Assumptions:" section
#    plus the final constraint line in front of it, so the bookmark stays
#    attached to the final line exactly as before.
$bookmarkPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$bookmarkPara.Range.InsertBefore(
    "Assumptions:" + [char]13 +
    "Admins will be inserted into the database through some method other than the application" + [char]13 +
    "Project names are unique")
